$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Row 11 / column B ("R40") is re-labelled "1". The target cell keeps its
# existing number format (General) but the new content must still be
# stored as TEXT, not a number, so a plain Value="1" assignment (which
# Excel auto-coerces to a numeric 1) is not enough. We stage the text in
# a scratch cell via a text-producing formula, copy it, and paste only
# the value into B11 - this carries the string type over without
# disturbing B11's existing style/formatting.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

